# Update the timestamp embedded in the test e-mail addresses on the
# "UsuariosRegistro" sheet (column C, rows 2-6) from 20251110_130229
# to 20251111_202811. Only this sheet's shared strings are touched,
# matching the source diff (the "LoginData" sheet keeps its own,
# untouched, older sample addresses).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UsuariosRegistro")

$oldStamp = "20251110_130229"
$newStamp = "20251111_202811"

$emailCol = 3  # Column C = "E-Mail"
for ($row = 2; $row -le 6; $row++) {
    $cell = $ws.Cells.Item($row, $emailCol)
    $current = $cell.Value2
    $cell.Value = $current -replace $oldStamp, $newStamp
}
